# Update LR-pair data with newly computed TPM-based values.
# The target cluster "ECs" is no longer a valid receiving cluster, so those
# rows are dropped; remaining rows (ECs/FAPs/MuSCs sending -> FAPs/MuSCs
# receiving) are recomputed with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (values as they will land in A2:T7).
$data = @(
    @("ECs",   "Adam10", "Epha3", "FAPs",  3, 1, 18.72828566666667, 56.184857,          0.37787211702967,    0.37787211702967,    3, 1, 26.097779, 78.29333700000001, 0.9922055808976035, 0.9922055808976036, 488.7666603775344,  4398.89994339781,  0.3749268233824309,  0.3749268233824309),
    @("ECs",   "Adam10", "Epha3", "MuSCs", 3, 1, 18.72828566666667, 56.184857,          0.37787211702967,    0.37787211702967,    3, 1, 0.205015,  0.6150450000000001, 0.007794419102396499, 0.007794419102396499, 3.839579485951667,  34.556215373565,  0.002945293647239065, 0.002945293647239065),
    @("FAPs",  "Adam10", "Epha3", "FAPs",  3, 1, 18.03569333333333, 54.10708,           0.3638979959652422,  0.3638979959652423,  3, 1, 26.097779, 78.29333700000001, 0.9922055808976035, 0.9922055808976036, 470.6915387251066,  4236.22384852596, 0.361061622474167,   0.3610616224741671),
    @("FAPs",  "Adam10", "Epha3", "MuSCs", 3, 1, 18.03569333333333, 54.10708,           0.3638979959652422,  0.3638979959652423,  3, 1, 0.205015,  0.6150450000000001, 0.007794419102396499, 0.007794419102396499, 3.697587668733334,  33.2782890186,    0.002836373491075288, 0.002836373491075289),
    @("MuSCs", "Adam10", "Epha3", "FAPs",  3, 1, 12.798518,          38.395554,         0.2582298870050878,  0.2582298870050878,  3, 1, 26.097779, 78.29333700000001, 0.9922055808976035, 0.9922055808976036, 334.0128942915221,  3006.116048623699, 0.2562171350410056,  0.2562171350410056),
    @("MuSCs", "Adam10", "Epha3", "MuSCs", 3, 1, 12.798518,          38.395554,         0.2582298870050878,  0.2582298870050878,  3, 1, 0.205015,  0.6150450000000001, 0.007794419102396499, 0.007794419102396499, 2.623888167770001,  23.61499350993001, 0.002012751964082145, 0.002012751964082145)
)

# Clear the old data body (rows 2-10) before writing the new, shorter table.
$ws.Range("A2:T10").Clear()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
